# Updates the cryptos price list for the Sun May  7 16:26:14 UTC 2023
# GitHub Actions refresh: new Price (column D) and Volume(1h) (column E)
# values for every coin row, plus a swap of the NEARProtocol/Quant rows
# (rows 49 and 50 exchange ranking position).
#
# Column D holds the price as plain text (values like "29.132.69" or
# "1.006" are NOT real numbers - they are pre-formatted strings from the
# scraped page). Setting NumberFormat to "@" (Text) before writing the
# value keeps Excel/the engine from re-interpreting those digit-and-dot
# strings as numeric literals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.111.18'
$ws.Range("E2").Value = '  +1.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.929.13'
$ws.Range("E3").Value = '  +2.41%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.57%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.48'
$ws.Range("E5").Value = '  +1.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  +0.57%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4613'
$ws.Range("E7").Value = '  +1.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3832'
$ws.Range("E8").Value = '  +1.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07784'
$ws.Range("E9").Value = '  +1.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9806'
$ws.Range("E10").Value = '  +2.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.71'
$ws.Range("E11").Value = '  +3.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.915.43'
$ws.Range("E12").Value = '  +1.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.703'
$ws.Range("E13").Value = '  +1.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.978'
$ws.Range("E14").Value = '  +0.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07088'
$ws.Range("E15").Value = '  +1.47%  '
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '84.51'
$ws.Range("E17").Value = '  +2.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009546'
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.77'
$ws.Range("E19").Value = '  +1.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.006'
$ws.Range("E20").Value = '  +0.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.146.40'
$ws.Range("E21").Value = '  +1.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.360'
$ws.Range("E22").Value = '  +1.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.98'
$ws.Range("E23").Value = '  +1.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.082'
$ws.Range("E24").Value = '  +0.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.10'
$ws.Range("E25").Value = '  +2.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.15'
$ws.Range("E26").Value = '  +1.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.689'
$ws.Range("E27").Value = '  +1.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '118.22'
$ws.Range("E28").Value = '  +1.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.847'
$ws.Range("E29").Value = '  +2.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09350'
$ws.Range("E30").Value = '  +1.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.8623'
$ws.Range("E31").Value = '  +2.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.136'
$ws.Range("E32").Value = '  +1.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.250'
$ws.Range("E33").Value = '  +0.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.032'
$ws.Range("E34").Value = '  +1.81%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05708'
$ws.Range("E35").Value = '  +1.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.162'
$ws.Range("E36").Value = '  +1.91%  '
$ws.Range("E37").Value = '  +0.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02053'
$ws.Range("E38").Value = '  +1.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.115'
$ws.Range("E39").Value = '  +16.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.542'
$ws.Range("E40").Value = '  +1.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5525'
$ws.Range("E41").Value = '  +1.23%  '
$ws.Range("E42").Value = '  +1.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.345'
$ws.Range("E43").Value = '  +2.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.000002830'
$ws.Range("E44").Value = '  -5.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.201'
$ws.Range("E45").Value = '  +7.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5218'
$ws.Range("E46").Value = '  +1.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.25'
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06941'
$ws.Range("E48").Value = '  +2.38%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '110.54'
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.774'
$ws.Range("E50").Value = '  +0.62%  '
$ws.Range("E51").Value = '  +0.58%  '
